# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recalculated and need to be
# overwritten with the new values below (rows 2-23 of Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 5
    10 = 0
    11 = 2
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 3
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 2
    22 = 1
    23 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
